# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Ajo" (Rosado variety) at rows 104-105,
# shifting the existing rows 104-122 down to rows 106-124.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 104 (twice, one row at a time)
$ws.Rows.Item(104).Insert()
$ws.Rows.Item(104).Insert()

# --- New row 104 ---
$ws.Range("A104").Value2 = 9
$ws.Range("B104").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C104").Value2 = "Metropolitana"
$ws.Range("D104").Value2 = 44505
$ws.Range("D104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E104").Value2 = 13
$ws.Range("F104").Value2 = 100112003
$ws.Range("G104").Value2 = "Ajo"
$ws.Range("H104").Value2 = "Rosado"
$ws.Range("I104").Value2 = "1a nueva(o)"
$ws.Range("J104").Value2 = 50
$ws.Range("K104").Value2 = 3400
$ws.Range("L104").Value2 = 3600
$ws.Range("M104").Value2 = 3500
$ws.Range("N104").Value2 = "$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O104").Value2 = "Provincia de Talagante"
$ws.Range("P104").Value2 = 175
$ws.Range("Q104").Value2 = 20
$ws.Range("R104").Value2 = "Hortaliza"

# --- New row 105 ---
$ws.Range("A105").Value2 = 9
$ws.Range("B105").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C105").Value2 = "Metropolitana"
$ws.Range("D105").Value2 = 44505
$ws.Range("D105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E105").Value2 = 13
$ws.Range("F105").Value2 = 100112003
$ws.Range("G105").Value2 = "Ajo"
$ws.Range("H105").Value2 = "Rosado"
$ws.Range("I105").Value2 = "Extra nueva (o)"
$ws.Range("J105").Value2 = 70
$ws.Range("K105").Value2 = 3800
$ws.Range("L105").Value2 = 4000
$ws.Range("M105").Value2 = 3900
$ws.Range("N105").Value2 = "$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O105").Value2 = "Provincia de Talagante"
$ws.Range("P105").Value2 = 195
$ws.Range("Q105").Value2 = 20
$ws.Range("R105").Value2 = "Hortaliza"
